# Edit slide 5's "Steps" text box (the second shape named "Title 1", a
# free-floating text box -- not the slide's title placeholder) to:
#   1. grow the shape's height (cy) slightly to fit the re-flowed text
#   2. update several bullet lines (the "State hypothesis" bullet is
#      removed and every following bullet's wording shifts up one slot),
#      add a trailing space after "Hypothesis :", merge the "uni"/"morgage"
#      typo-fix runs into single clean runs, and fix "morgage" -> "mortgage"

$pres = $ppt.ActivePresentation
$slide = $pres.Slides.Item(5)
$shape = $slide.Shapes.Item(4)

# --- 1. Resize the shape (cy 4857749 -> 5052059; cx/off unchanged) ---
# Shape.Top/Left/Width/Height are expressed in points (EMU / 12700), so
# convert the target EMU value rather than assigning the raw EMU number.
# Only Height actually changes -- leave Top/Left/Width alone so their
# already-correct EMU values are not disturbed by a points round-trip.
$shape.Height = 5052059 / 12700

# --- 2. Rewrite the bulleted text, paragraph by paragraph -----------
# Paragraphs are 1-indexed character runs inside one big TextRange,
# separated by `r`. Working from the LAST paragraph back to the FIRST
# keeps every not-yet-touched offset valid, since only edits to later
# text can shift earlier offsets.

$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 13 (index from 1): "What if the wages ... morgage rate ..."
# -> single run, "morgage" fixed to "mortgage"
$tr.Characters(461, 177).Text = "What if the wages continue to rise at their current rate and the mortgage rate (30 year fixed) holds steady at 4.6 % ...Can you afford to buy a median price home in the future???"

# Paragraph 10: "Use the coefficients and y-intercept to predict future outcomes" -> removed (paragraph becomes empty)
$tr.Characters(360, 63).Text = ""

# Paragraph 9: "Run multi-variate linear regression model for all 3 predictor variable" -> "Use the coefficients and y-intercept to predict future outcomes"
$tr.Characters(289, 70).Text = "Use the coefficients and y-intercept to predict future outcomes"

# Paragraph 8: "run " + "uni" + "-variate linear regression model for each predictor variable" -> "Run multi-variate linear regression model for all 3 predictor variable"
$tr.Characters(221, 67).Text = "Run multi-variate linear regression model for all 3 predictor variable"

# Paragraph 7: "Short-list 3 most correlated predictor variables " -> "run uni-variate linear regression model for each predictor variable"
$tr.Characters(171, 49).Text = "run uni-variate linear regression model for each predictor variable"

# Paragraph 6: "Check initial correlation between dependent and potential predictor variables" -> "Short-list 3 most correlated predictor variables "
$tr.Characters(93, 77).Text = "Short-list 3 most correlated predictor variables "

# Paragraph 5: "Identify dependent and potential predictor variables" -> "Check initial correlation between dependent and potential predictor variables"
$tr.Characters(40, 52).Text = "Check initial correlation between dependent and potential predictor variables"

# Paragraph 4: "State hypothesis" -> "Identify dependent and potential predictor variables"
$tr.Characters(23, 16).Text = "Identify dependent and potential predictor variables"

# Paragraph 1: "Hypothesis :" -> "Hypothesis : " (trailing space added)
$tr.Characters(1, 12).Text = "Hypothesis : "
